{"js": "// Replace the 25 \"three-digit \u00d7 one-digit\" equations in the table with\n// their new values, each as a single self-contained search/replace so\n// the run/paragraph formatting (font, size, justification) is preserved.\nconst replacements = [\n  [\"486\u00d75=2430\", \"824\u00d74=3296\"],\n  [\"837\u00d77=5859\", \"257\u00d78=2056\"],\n  [\"486\u00d79=4374\", \"114\u00d73=342\"],\n  [\"750\u00d74=3000\", \"703\u00d79=6327\"],\n  [\"981\u00d78=7848\", \"450\u00d72=900\"],\n  [\"478\u00d73=1434\", \"842\u00d75=4210\"],\n  [\"609\u00d77=4263\", \"878\u00d74=3512\"],\n  [\"702\u00d72=1404\", \"210\u00d74=840\"],\n  [\"577\u00d75=2885\", \"196\u00d79=1764\"],\n  [\"719\u00d78=5752\", \"411\u00d78=3288\"],\n  [\"612\u00d73=1836\", \"330\u00d78=2640\"],\n  [\"137\u00d79=1233\", \"750\u00d76=4500\"],\n  [\"141\u00d79=1269\", \"531\u00d75=2655\"],\n  [\"371\u00d72=742\", \"415\u00d74=1660\"],\n  [\"308\u00d79=2772\", \"433\u00d73=1299\"],\n  [\"714\u00d78=5712\", \"646\u00d76=3876\"],\n  [\"971\u00d76=5826\", \"400\u00d79=3600\"],\n  [\"256\u00d75=1280\", \"254\u00d78=2032\"],\n  [\"640\u00d75=3200\", \"651\u00d74=2604\"],\n  [\"206\u00d74=824\", \"636\u00d79=5724\"],\n  [\"416\u00d72=832\", \"652\u00d74=2608\"],\n  [\"720\u00d75=3600\", \"986\u00d72=1972\"],\n  [\"399\u00d74=1596\", \"107\u00d78=856\"],\n  [\"830\u00d76=4980\", \"614\u00d72=1228\"],\n  [\"730\u00d73=2190\", \"841\u00d76=5046\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"three-digit \u00d7 one-digit\" equations in the table with\n# their new values. Each pair is applied with Find/Replace scoped to the\n# whole document content, matching case exactly, with no wildcards, so\n# only the intended full-cell text is touched and run formatting (font,\n# size, paragraph alignment) is left untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"486\u00d75=2430\", \"824\u00d74=3296\"),\n  @(\"837\u00d77=5859\", \"257\u00d78=2056\"),\n  @(\"486\u00d79=4374\", \"114\u00d73=342\"),\n  @(\"750\u00d74=3000\", \"703\u00d79=6327\"),\n  @(\"981\u00d78=7848\", \"450\u00d72=900\"),\n  @(\"478\u00d73=1434\", \"842\u00d75=4210\"),\n  @(\"609\u00d77=4263\", \"878\u00d74=3512\"),\n  @(\"702\u00d72=1404\", \"210\u00d74=840\"),\n  @(\"577\u00d75=2885\", \"196\u00d79=1764\"),\n  @(\"719\u00d78=5752\", \"411\u00d78=3288\"),\n  @(\"612\u00d73=1836\", \"330\u00d78=2640\"),\n  @(\"137\u00d79=1233\", \"750\u00d76=4500\"),\n  @(\"141\u00d79=1269\", \"531\u00d75=2655\"),\n  @(\"371\u00d72=742\", \"415\u00d74=1660\"),\n  @(\"308\u00d79=2772\", \"433\u00d73=1299\"),\n  @(\"714\u00d78=5712\", \"646\u00d76=3876\"),\n  @(\"971\u00d76=5826\", \"400\u00d79=3600\"),\n  @(\"256\u00d75=1280\", \"254\u00d78=2032\"),\n  @(\"640\u00d75=3200\", \"651\u00d74=2604\"),\n  @(\"206\u00d74=824\", \"636\u00d79=5724\"),\n  @(\"416\u00d72=832\", \"652\u00d74=2608\"),\n  @(\"720\u00d75=3600\", \"986\u00d72=1972\"),\n  @(\"399\u00d74=1596\", \"107\u00d78=856\"),\n  @(\"830\u00d76=4980\", \"614\u00d72=1228\"),\n  @(\"730\u00d73=2190\", \"841\u00d76=5046\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
